# Scheduled-runner market data refresh: push updated currentAveragePrice /
# LevePrice / LeveProfit figures into the per-job leve tables (one table per
# worksheet/job: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 228.25
$ws.Range("I5").Value = 228.25
$ws.Range("K5").Value = 228.25
$ws.Range("M5").Value = -113.25

$ws.Range("H17").Value = 866.25
$ws.Range("I17").Value = 400
$ws.Range("J17").Value = 918.05554
$ws.Range("K17").Value = 1200
$ws.Range("L17").Value = 2754.16662
$ws.Range("M17").Value = -1032
$ws.Range("N17").Value = -3090.16662

$ws.Range("H19").Value = 2369.7083
$ws.Range("I19").Value = 1754.625
$ws.Range("J19").Value = 2677.25
$ws.Range("K19").Value = 1754.625
$ws.Range("L19").Value = 2677.25
$ws.Range("M19").Value = -1579.625
$ws.Range("N19").Value = -3027.25

$ws.Range("H70").Value = 144249.86
$ws.Range("I70").Value = 1624.8334
$ws.Range("K70").Value = 4874.5002
$ws.Range("M70").Value = -4604.5002

$ws.Range("H73").Value = 144249.86
$ws.Range("I73").Value = 1624.8334
$ws.Range("K73").Value = 4874.5002
$ws.Range("M73").Value = -3938.5002

$ws.Range("H112").Value = 3195.8333
$ws.Range("J112").Value = 3269.6086
$ws.Range("L112").Value = 9808.825800000001
$ws.Range("N112").Value = -12024.8258

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8846.484
$ws.Range("J45").Value = 4444.6665
$ws.Range("L45").Value = 4444.6665
$ws.Range("N45").Value = -5198.6665

$ws.Range("H64").Value = 45000
$ws.Range("J64").Value = 45000
$ws.Range("L64").Value = 45000
$ws.Range("N64").Value = -45496

$ws.Range("H67").Value = 45000
$ws.Range("J67").Value = 45000
$ws.Range("L67").Value = 45000
$ws.Range("N67").Value = -46716

$ws.Range("H131").Value = 84715
$ws.Range("J131").Value = 84715
$ws.Range("L131").Value = 84715
$ws.Range("N131").Value = -94795

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1900.4348
$ws.Range("I20").Value = 2134.111
$ws.Range("K20").Value = 2134.111
$ws.Range("M20").Value = -1887.111

$ws.Range("H134").Value = 1610.3846
$ws.Range("I134").Value = 1540.091
$ws.Range("J134").Value = 1997
$ws.Range("K134").Value = 4620.272999999999
$ws.Range("L134").Value = 5991
$ws.Range("M134").Value = -2085.272999999999
$ws.Range("N134").Value = -11061

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1210.3334
$ws.Range("I107").Value = 1003.5789
$ws.Range("K107").Value = 1003.5789
$ws.Range("M107").Value = 916.4211

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 963.2727
$ws.Range("J107").Value = 1339.6
$ws.Range("L107").Value = 4018.8
$ws.Range("N107").Value = -7858.799999999999

$ws.Range("H111").Value = 8140.7144
$ws.Range("I111").Value = 7397.2
$ws.Range("K111").Value = 22191.6
$ws.Range("M111").Value = -19124.6

$ws.Range("H117").Value = 3184
$ws.Range("I117").Value = 500
$ws.Range("J117").Value = 3482.2222
$ws.Range("K117").Value = 1500
$ws.Range("L117").Value = 10446.6666
$ws.Range("M117").Value = 1942
$ws.Range("N117").Value = -17330.6666

$ws.Range("H120").Value = 15645.429
$ws.Range("I120").Value = 13253.167
$ws.Range("K120").Value = 39759.501
$ws.Range("M120").Value = -34921.501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 21666.666
$ws.Range("J21").Value = 21666.666
$ws.Range("L21").Value = 21666.666
$ws.Range("N21").Value = -22012.666

$ws.Range("H30").Value = 21666.666
$ws.Range("J30").Value = 21666.666
$ws.Range("L30").Value = 21666.666
$ws.Range("N30").Value = -21876.666

$ws.Range("H35").Value = 4506.5
$ws.Range("I35").Value = 4506.5
$ws.Range("K35").Value = 4506.5
$ws.Range("M35").Value = -4208.5

$ws.Range("H70").Value = 5597.2
$ws.Range("I70").Value = 6000
$ws.Range("J70").Value = 5328.6665
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 5328.6665
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -5868.6665

$ws.Range("H73").Value = 5597.2
$ws.Range("I73").Value = 6000
$ws.Range("J73").Value = 5328.6665
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 5328.6665
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -7200.6665

$ws.Range("H80").Value = 4194.9414
$ws.Range("J80").Value = 5923.1
$ws.Range("L80").Value = 5923.1
$ws.Range("N80").Value = -7919.1

$ws.Range("H83").Value = 4194.9414
$ws.Range("J83").Value = 5923.1
$ws.Range("L83").Value = 29615.5
$ws.Range("N83").Value = -39599.5

$ws.Range("H102").Value = 3645.2222
$ws.Range("I102").Value = 3686.7144
$ws.Range("K102").Value = 3686.7144
$ws.Range("M102").Value = -2064.7144

$ws.Range("H122").Value = 9549.5
$ws.Range("I122").Value = 5389.3
$ws.Range("K122").Value = 16167.9
$ws.Range("M122").Value = -13717.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H68").Value = 2748.75
$ws.Range("I68").Value = 3499
$ws.Range("J68").Value = 1998.5
$ws.Range("K68").Value = 3499
$ws.Range("L68").Value = 1998.5
$ws.Range("M68").Value = -2750
$ws.Range("N68").Value = -3496.5

$ws.Range("H71").Value = 2748.75
$ws.Range("I71").Value = 3499
$ws.Range("J71").Value = 1998.5
$ws.Range("K71").Value = 17495
$ws.Range("L71").Value = 9992.5
$ws.Range("M71").Value = -13751
$ws.Range("N71").Value = -17480.5

$ws.Range("H82").Value = 993
$ws.Range("I82").Value = 848.2857
$ws.Range("K82").Value = 848.2857
$ws.Range("M82").Value = -487.2857

$ws.Range("H85").Value = 993
$ws.Range("I85").Value = 848.2857
$ws.Range("K85").Value = 848.2857
$ws.Range("M85").Value = 399.7143

$ws.Range("H122").Value = 4466
$ws.Range("I122").Value = 4599
$ws.Range("K122").Value = 13797
$ws.Range("M122").Value = -11347

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7826
$ws.Range("I62").Value = 5998
$ws.Range("K62").Value = 5998
$ws.Range("M62").Value = -5374

$ws.Range("H65").Value = 7826
$ws.Range("I65").Value = 5998
$ws.Range("K65").Value = 29990
$ws.Range("M65").Value = -26870

$ws.Range("H122").Value = 8354.9
$ws.Range("I122").Value = 7274.8335
$ws.Range("K122").Value = 21824.5005
$ws.Range("M122").Value = -19374.5005

$ws.Range("H136").Value = 1695.5416
$ws.Range("J136").Value = 2397.5
$ws.Range("L136").Value = 7192.5
$ws.Range("N136").Value = -12292.5
